$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @("SE", "SE", "FICT", "SE", "SE", "SE", "SE", "SE")

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 8 + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
